$wb = $excel.ActiveWorkbook

# Both the "展览" (Exhibition) sheet and the "全部类型" (All types) sheet
# contain the same underlying data and need the same F-column ("想去人数")
# updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1948
    $ws.Range("F4").Value = 1218
    $ws.Range("F5").Value = 1376
    $ws.Range("F6").Value = 52
    $ws.Range("F7").Value = 6081
}
